# Re-upload of the roster reorders two player rows:
#   - "Derrick White" (previously the very last player row, row 19) moves
#     up to sit right after "Quentin Grimes" (new row 5), pushing
#     "Harrison Barnes" .. "Jared McCain" down by one row each.
#   - "Cade Cunningham" (previously row 13, right after "Jared McCain")
#     moves down to become the new very last row (row 19), with
#     "Harrison Barnes" .. "Jared McCain" now occupying rows 6:13 and
#     "Damian Lillard" .. "Brandon Ingram" staying put at rows 14:18.
# Net effect: the two rows swap ends of the table while everything
# between them shifts by one row to close/open the gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove "Cade Cunningham" (row 13). Rows 14:19 shift up to 13:18.
$ws.Rows("13:13").Delete()

# 2) Remove "Derrick White", now the last row (18). Nothing above shifts.
$ws.Rows("18:18").Delete()

# 3) Insert a fresh blank row at 5 for "Derrick White"; rows 5:17 shift
#    down to 6:18, restoring 18 data rows below the header.
$ws.Rows("5:5").Insert()
$ws.Range("A5").Value = "Derrick White"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Boston Celtics"

# 4) Append "Cade Cunningham" back as the new last row (19).
$ws.Range("A19").Value = "Cade Cunningham"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Detroit Pistons"
